$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "------"
$ws.Range("B2").Value = "----"
$ws.Range("C2").Value = "-----"
$ws.Range("D2").Value = "-----------"
$ws.Range("E2").Value = "----------"
$ws.Range("F2").Value = "--------"
$ws.Range("G2").Value = "--------"
$ws.Range("H2").Value = "--------"
$ws.Range("I2").Value = "-----"
$ws.Range("J2").Value = "-----"
$ws.Range("K2").Value = "-----"
$ws.Range("L2").Value = "-----"
$ws.Range("M2").Value = "------"
$ws.Range("N2").Value = "-----"
$ws.Range("O2").Value = "-----"
$ws.Range("P2").Value = "-----"

$ws.Range("A3").Value = 6380
$ws.Range("B3").Value = 4.31
$ws.Range("C3").Value = 0.2
$ws.Range("D3").Value = "4-parameter"
$ws.Range("E3").Value = "STIS.G430L"
$ws.Range("F3").Value = 0.45
$ws.Range("G3").Value = 0.4807
$ws.Range("H3").Value = 0.5117
$ws.Range("I3").Value = 0.229
$ws.Range("J3").Value = 0.051
$ws.Range("K3").Value = 1.118
$ws.Range("L3").Value = 0.121
$ws.Range("M3").Value = -0.755
$ws.Range("N3").Value = 0.122
$ws.Range("O3").Value = 0.221
$ws.Range("P3").Value = 0.044

$ws.Range("A4").Value = 6380
$ws.Range("B4").Value = 4.31
$ws.Range("C4").Value = 0.2
$ws.Range("D4").Value = "4-parameter"
$ws.Range("E4").Value = "STIS.G430L"
$ws.Range("F4").Value = 0.512
$ws.Range("G4").Value = 0.5425
$ws.Range("H4").Value = 0.573
$ws.Range("I4").Value = 0.157
$ws.Range("J4").Value = 0.043
$ws.Range("K4").Value = 1.287
$ws.Range("L4").Value = 0.102
$ws.Range("M4").Value = -1.023
$ws.Range("N4").Value = 0.103
$ws.Range("O4").Value = 0.324
$ws.Range("P4").Value = 0.037

$ws.Range("A5").Value = 6380
$ws.Range("B5").Value = 4.31
$ws.Range("C5").Value = 0.2
$ws.Range("D5").Value = "4-parameter"
$ws.Range("E5").Value = "STIS.G750L"
$ws.Range("F5").Value = 0.525
$ws.Range("G5").Value = 0.575
$ws.Range("H5").Value = 0.626
$ws.Range("I5").Value = 0.125
$ws.Range("J5").Value = 0.038
$ws.Range("K5").Value = 1.395
$ws.Range("L5").Value = 0.09
$ws.Range("M5").Value = -1.196
$ws.Range("N5").Value = 0.091
$ws.Range("O5").Value = 0.39
$ws.Range("P5").Value = 0.033

$ws.Range("A6").Value = 6380
$ws.Range("B6").Value = 4.31
$ws.Range("C6").Value = 0.2
$ws.Range("D6").Value = "4-parameter"
$ws.Range("E6").Value = "STIS.G750L"
$ws.Range("F6").Value = 0.626
$ws.Range("G6").Value = 0.677
$ws.Range("H6").Value = 0.727
$ws.Range("I6").Value = 0.074
$ws.Range("J6").Value = 0.027
$ws.Range("K6").Value = 1.469
$ws.Range("L6").Value = 0.063
$ws.Range("M6").Value = -1.35
$ws.Range("N6").Value = 0.064
$ws.Range("O6").Value = 0.444
$ws.Range("P6").Value = 0.023

$ws.Range("A7").Value = 6380
$ws.Range("B7").Value = 4.31
$ws.Range("C7").Value = 0.2
$ws.Range("D7").Value = "4-parameter"
$ws.Range("E7").Value = "STIS.G750L"
$ws.Range("F7").Value = 0.728
$ws.Range("G7").Value = 0.7783
$ws.Range("H7").Value = 0.829
$ws.Range("I7").Value = 0.058
$ws.Range("J7").Value = 0.018
$ws.Range("K7").Value = 1.414
$ws.Range("L7").Value = 0.042
$ws.Range("M7").Value = -1.338
$ws.Range("N7").Value = 0.043
$ws.Range("O7").Value = 0.445
$ws.Range("P7").Value = 0.016

$ws.Range("A8").Value = 6380
$ws.Range("B8").Value = 4.31
$ws.Range("C8").Value = 0.2
$ws.Range("D8").Value = "4-parameter"
$ws.Range("E8").Value = "STIS.G750L"
$ws.Range("F8").Value = 0.829
$ws.Range("G8").Value = 0.88
$ws.Range("H8").Value = 0.93
$ws.Range("I8").Value = 0.146
$ws.Range("J8").Value = 0.022
$ws.Range("K8").Value = 1.118
$ws.Range("L8").Value = 0.052
$ws.Range("M8").Value = -1.069
$ws.Range("N8").Value = 0.052
$ws.Range("O8").Value = 0.345
$ws.Range("P8").Value = 0.019

$ws.Range("A9").Value = 6380
$ws.Range("B9").Value = 4.31
$ws.Range("C9").Value = 0.2
$ws.Range("D9").Value = "4-parameter"
$ws.Range("E9").Value = "STIS.G750L"
$ws.Range("F9").Value = 0.931
$ws.Range("G9").Value = 0.9814
$ws.Range("H9").Value = 1.032
$ws.Range("I9").Value = 0.149
$ws.Range("J9").Value = 0.016
$ws.Range("K9").Value = 1.085
$ws.Range("L9").Value = 0.039
$ws.Range("M9").Value = -1.076
$ws.Range("N9").Value = 0.039
$ws.Range("O9").Value = 0.363
$ws.Range("P9").Value = 0.014

$ws.Range("A10").Value = 6380
$ws.Range("B10").Value = 4.31
$ws.Range("C10").Value = 0.2
$ws.Range("D10").Value = "4-parameter"
$ws.Range("E10").Value = "HST/WFC3_IR.G141"
$ws.Range("F10").Value = 1.04
$ws.Range("G10").Value = 1.098
$ws.Range("H10").Value = 1.157
$ws.Range("I10").Value = 0.241
$ws.Range("J10").Value = 0.02
$ws.Range("K10").Value = 0.878
$ws.Range("L10").Value = 0.048
$ws.Range("M10").Value = -0.949
$ws.Range("N10").Value = 0.048
$ws.Range("O10").Value = 0.333
$ws.Range("P10").Value = 0.018

$ws.Range("A11").Value = 6380
$ws.Range("B11").Value = 4.31
$ws.Range("C11").Value = 0.2
$ws.Range("D11").Value = "4-parameter"
$ws.Range("E11").Value = "HST/WFC3_IR.G141"
$ws.Range("F11").Value = 1.158
$ws.Range("G11").Value = 1.222
$ws.Range("H11").Value = 1.288
$ws.Range("I11").Value = 0.357
$ws.Range("J11").Value = 0.025
$ws.Range("K11").Value = 0.663
$ws.Range("L11").Value = 0.058
$ws.Range("M11").Value = -0.813
$ws.Range("N11").Value = 0.058
$ws.Range("O11").Value = 0.301
$ws.Range("P11").Value = 0.021

$ws.Range("A12").Value = 6380
$ws.Range("B12").Value = 4.31
$ws.Range("C12").Value = 0.2
$ws.Range("D12").Value = "4-parameter"
$ws.Range("E12").Value = "HST/WFC3_IR.G141"
$ws.Range("F12").Value = 1.288
$ws.Range("G12").Value = 1.359
$ws.Range("H12").Value = 1.434
$ws.Range("I12").Value = 0.573
$ws.Range("J12").Value = 0.023
$ws.Range("K12").Value = 0.224
$ws.Range("L12").Value = 0.055
$ws.Range("M12").Value = -0.478
$ws.Range("N12").Value = 0.056
$ws.Range("O12").Value = 0.204
$ws.Range("P12").Value = 0.02

$ws.Range("A13").Value = 6380
$ws.Range("B13").Value = 4.31
$ws.Range("C13").Value = 0.2
$ws.Range("D13").Value = "4-parameter"
$ws.Range("E13").Value = "HST/WFC3_IR.G141"
$ws.Range("F13").Value = 1.434
$ws.Range("G13").Value = 1.513
$ws.Range("H13").Value = 1.595
$ws.Range("I13").Value = 0.899
$ws.Range("J13").Value = 0.014
$ws.Range("K13").Value = -0.536
$ws.Range("L13").Value = 0.032
$ws.Range("M13").Value = 0.196
$ws.Range("N13").Value = 0.032
$ws.Range("O13").Value = -0.019
$ws.Range("P13").Value = 0.012

$ws.Range("A14").Value = 6380
$ws.Range("B14").Value = 4.31
$ws.Range("C14").Value = 0.2
$ws.Range("D14").Value = "4-parameter"
$ws.Range("E14").Value = "HST/WFC3_IR.G141"
$ws.Range("F14").Value = 1.595
$ws.Range("G14").Value = 1.683
$ws.Range("H14").Value = 1.774
$ws.Range("I14").Value = 1.013
$ws.Range("J14").Value = 0.007
$ws.Range("K14").Value = -0.85
$ws.Range("L14").Value = 0.017
$ws.Range("M14").Value = 0.469
$ws.Range("N14").Value = 0.017
$ws.Range("O14").Value = -0.108
$ws.Range("P14").Value = 0.006

$ws.Range("AB22").Select() | Out-Null
